$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.061.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.21%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.088.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.27%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.18%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.93%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.086.93'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.32%  '

# Row 9
$ws.Range('E9').Value = '  -0.41%  '

# Row 10
$ws.Range('E10').Value = '  -0.53%  '

# Row 11
$ws.Range('E11').Value = '  +0.94%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.465'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.37%  '

# Row 13
$ws.Range('E13').Value = '  -0.27%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.18%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.604.63'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.34%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.995.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.28%  '

# Row 18
$ws.Range('E18').Value = '  -0.42%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.089.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.23%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.37%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '487.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.52%  '

# Row 23
$ws.Range('E23').Value = '  -0.79%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.29'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.26%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.41%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.74%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.04%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.55%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.38%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.09%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.50%  '

# Row 33
$ws.Range('E33').Value = '  -0.17%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0939'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.09%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.09%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '47.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.86%  '

# Row 37
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.942'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.16%  '

# Row 38
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.56%  '

# Row 39
$ws.Range('E39').Value = '  +2.82%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.56%  '

# Row 41
$ws.Range('E41').Value = '  +1.17%  '

# Row 42
$ws.Range('E42').Value = '  -0.30%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.72'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.33%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.37%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.785.69'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.47%  '

# Row 46
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '367.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.46%  '

# Row 47
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0343'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '

# Row 48
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.42%  '

# Row 49
$ws.Range('E49').Value = '  +0.04%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.15%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.40%  '
